# Batch runner is added
# Mark the SuiteType of test cases 49-163 (rows 52-167, skipping the blank
# separator rows) on the "SuiteDetails" sheet as "Regression" instead of
# "SmokeSuite", and move the view/selection down to the new bottom of the
# list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SuiteDetails")

# Column I ("SuiteType") rows 52-167, skipping the 4 blank separator rows
# (94, 117, 136, 159) that delimit sub-blocks of test cases.
$ws.Range("I52:I93").Value   = "Regression"
$ws.Range("I95:I116").Value  = "Regression"
$ws.Range("I118:I135").Value = "Regression"
$ws.Range("I137:I158").Value = "Regression"
$ws.Range("I160:I167").Value = "Regression"

# Scroll the frozen view down to the newly edited rows and move the active
# selection to just past the last data row.
$win = $excel.ActiveWindow
$win.ScrollRow = 153
$win.ScrollColumn = 2
$ws.Range("G168").Select()
